$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# First-page footer (footer1.xml, wp:docPr id="3") — Pearson logo.
$ftrFirst = $sec.Footers.Item(2)
$ftrFirst.Range.InlineShapes.Item(1).Name = "image1.png"

# Default footer (footer2.xml, wp:docPr id="2") — Pearson logo.
$ftrDefault = $sec.Footers.Item(1)
$ftrDefault.Range.InlineShapes.Item(1).Name = "image1.png"

# First-page header (header1.xml, wp:docPr id="1") — BTEC logo.
$hdrFirst = $sec.Headers.Item(2)
$hdrFirst.Range.InlineShapes.Item(1).Name = "image2.jpg"
